$p = $ppt.ActivePresentation

# Locate the slide/shape that contains the "def max(x, y): ... return x ... return y"
# code sample (Functions: Example slide).
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $t = $shp.TextFrame.TextRange.Text
                if ($t -like "*def max(x, y):*return x*return y*") {
                    $targetSlide = $sl
                    $targetShape = $shp
                }
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# --- 1) Split "def max(x, y):" into "def " / "print_max" / "(x, y):" ---
$defPara = $tr.Paragraphs(1)
$maxWord = $defPara.Find("max")
$maxWord.Text = "print_max"

# --- 2) "        return x" -> "        print( x )" ---
$retXPara = $tr.Paragraphs(3)
$retXFull = $tr.Characters($retXPara.Start, $retXPara.Length)
$retXFull.Text = "        print( x )"

# --- 3) "        return y" -> "        print ( y )" ---
$retYPara = $tr.Paragraphs(5)
$retYFull = $tr.Characters($retYPara.Start, $retYPara.Length)
$retYFull.Text = "        print ( y )"

# --- 4) Add a small empty textbox on the same slide (leftover artifact) ---
$emu = 12700
$tb = $targetSlide.Shapes.AddTextbox(1, 8172450 / $emu, 3814763 / $emu, 184731 / $emu, 369332 / $emu)
$tb.Fill.Visible = 0
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1

# Reorder it so it sits right after the code sample shape (matches the
# original author's shape order, ahead of the two connector shapes).
$tb.ZOrder(3)
$tb.ZOrder(3)
